$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the last existing header cell (G1) onto the
# new H1 header cell so it matches the other header cells' formatting
# (bold font, border, centered alignment).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Set the new header label and its value in row 2
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
